# Translation key clean-up ("Initial check-in of translations changes.")
#
# The "survey" sheet's prompt-level display header and the "settings"
# sheet's form-level display header are renamed to their more specific,
# namespaced equivalents:
#   display.text   -> display.prompt.text   (survey!C1)
#   display.title  -> display.title.text    (settings!C1)
#
# Order matters for shared-string allocation: settings is touched first so
# "display.title.text" gets the lower new string index, matching how the
# workbook was actually re-saved.

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("C1").Value = "display.title.text"

$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Range("C1").Value = "display.prompt.text"

# Final state of the workbook had "survey" as the active tab with C16
# selected there, and the cursor left on settings!C2 before moving away -
# reproduce the same view/selection bookkeeping.
$wsSurvey.Activate() | Out-Null
$wsSurvey.Range("C16").Select() | Out-Null

$wsSettings.Activate() | Out-Null
$wsSettings.Range("C2").Select() | Out-Null

$wsSurvey.Activate() | Out-Null
